$d = $word.ActiveDocument

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:widowControl/>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:spacing w:before="150" w:after="150"/>
    <w:jc w:val="left"/>
    <w:rPr>
      <w:rFonts w:ascii="微软雅黑" w:eastAsia="微软雅黑" w:hAnsi="微软雅黑" w:cs="宋体"/>
      <w:color w:val="333333"/>
      <w:kern w:val="0"/>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="微软雅黑" w:eastAsia="微软雅黑" w:hAnsi="微软雅黑" w:cs="宋体" w:hint="eastAsia"/>
      <w:color w:val="333333"/>
      <w:kern w:val="0"/>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t>在settings配置文件中有一个ITEM_PIPELINES的配置参数 :</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:widowControl/>
    <w:shd w:val="clear" w:color="auto" w:fill="F5F5F5"/>
    <w:tabs>
      <w:tab w:val="left" w:pos="916"/>
      <w:tab w:val="left" w:pos="1832"/>
      <w:tab w:val="left" w:pos="2748"/>
      <w:tab w:val="left" w:pos="3664"/>
      <w:tab w:val="left" w:pos="4580"/>
      <w:tab w:val="left" w:pos="5496"/>
      <w:tab w:val="left" w:pos="6412"/>
      <w:tab w:val="left" w:pos="7328"/>
      <w:tab w:val="left" w:pos="8244"/>
      <w:tab w:val="left" w:pos="9160"/>
      <w:tab w:val="left" w:pos="10076"/>
      <w:tab w:val="left" w:pos="10992"/>
      <w:tab w:val="left" w:pos="11908"/>
      <w:tab w:val="left" w:pos="12824"/>
      <w:tab w:val="left" w:pos="13740"/>
      <w:tab w:val="left" w:pos="14656"/>
    </w:tabs>
    <w:jc w:val="left"/>
    <w:rPr>
      <w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体" w:cs="宋体"/>
      <w:color w:val="000000"/>
      <w:kern w:val="0"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体" w:cs="宋体"/>
      <w:color w:val="000000"/>
      <w:kern w:val="0"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>ITEM_PIPELINES = {</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:widowControl/>
    <w:shd w:val="clear" w:color="auto" w:fill="F5F5F5"/>
    <w:tabs>
      <w:tab w:val="left" w:pos="916"/>
      <w:tab w:val="left" w:pos="1832"/>
      <w:tab w:val="left" w:pos="2748"/>
      <w:tab w:val="left" w:pos="3664"/>
      <w:tab w:val="left" w:pos="4580"/>
      <w:tab w:val="left" w:pos="5496"/>
      <w:tab w:val="left" w:pos="6412"/>
      <w:tab w:val="left" w:pos="7328"/>
      <w:tab w:val="left" w:pos="8244"/>
      <w:tab w:val="left" w:pos="9160"/>
      <w:tab w:val="left" w:pos="10076"/>
      <w:tab w:val="left" w:pos="10992"/>
      <w:tab w:val="left" w:pos="11908"/>
      <w:tab w:val="left" w:pos="12824"/>
      <w:tab w:val="left" w:pos="13740"/>
      <w:tab w:val="left" w:pos="14656"/>
    </w:tabs>
    <w:jc w:val="left"/>
    <w:rPr>
      <w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体" w:cs="宋体"/>
      <w:color w:val="000000"/>
      <w:kern w:val="0"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体" w:cs="宋体"/>
      <w:color w:val="000000"/>
      <w:kern w:val="0"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">    </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体" w:cs="宋体"/>
      <w:color w:val="800000"/>
      <w:kern w:val="0"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>'</w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体" w:cs="宋体"/>
      <w:color w:val="800000"/>
      <w:kern w:val="0"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>myproject.pipelines</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体" w:cs="宋体"/>
      <w:color w:val="800000"/>
      <w:kern w:val="0"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>.PricePipeline</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体" w:cs="宋体"/>
      <w:color w:val="800000"/>
      <w:kern w:val="0"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>'</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体" w:cs="宋体"/>
      <w:color w:val="000000"/>
      <w:kern w:val="0"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>: 300,</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:widowControl/>
    <w:shd w:val="clear" w:color="auto" w:fill="F5F5F5"/>
    <w:tabs>
      <w:tab w:val="left" w:pos="916"/>
      <w:tab w:val="left" w:pos="1832"/>
      <w:tab w:val="left" w:pos="2748"/>
      <w:tab w:val="left" w:pos="3664"/>
      <w:tab w:val="left" w:pos="4580"/>
      <w:tab w:val="left" w:pos="5496"/>
      <w:tab w:val="left" w:pos="6412"/>
      <w:tab w:val="left" w:pos="7328"/>
      <w:tab w:val="left" w:pos="8244"/>
      <w:tab w:val="left" w:pos="9160"/>
      <w:tab w:val="left" w:pos="10076"/>
      <w:tab w:val="left" w:pos="10992"/>
      <w:tab w:val="left" w:pos="11908"/>
      <w:tab w:val="left" w:pos="12824"/>
      <w:tab w:val="left" w:pos="13740"/>
      <w:tab w:val="left" w:pos="14656"/>
    </w:tabs>
    <w:jc w:val="left"/>
    <w:rPr>
      <w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体" w:cs="宋体"/>
      <w:color w:val="000000"/>
      <w:kern w:val="0"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体" w:cs="宋体"/>
      <w:color w:val="000000"/>
      <w:kern w:val="0"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">    </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体" w:cs="宋体"/>
      <w:color w:val="800000"/>
      <w:kern w:val="0"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>'</w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体" w:cs="宋体"/>
      <w:color w:val="800000"/>
      <w:kern w:val="0"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>myproject.pipelines</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体" w:cs="宋体"/>
      <w:color w:val="800000"/>
      <w:kern w:val="0"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>.JsonWriterPipeline</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体" w:cs="宋体"/>
      <w:color w:val="800000"/>
      <w:kern w:val="0"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>'</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体" w:cs="宋体"/>
      <w:color w:val="000000"/>
      <w:kern w:val="0"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>: 800,</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:widowControl/>
    <w:shd w:val="clear" w:color="auto" w:fill="F5F5F5"/>
    <w:tabs>
      <w:tab w:val="left" w:pos="916"/>
      <w:tab w:val="left" w:pos="1832"/>
      <w:tab w:val="left" w:pos="2748"/>
      <w:tab w:val="left" w:pos="3664"/>
      <w:tab w:val="left" w:pos="4580"/>
      <w:tab w:val="left" w:pos="5496"/>
      <w:tab w:val="left" w:pos="6412"/>
      <w:tab w:val="left" w:pos="7328"/>
      <w:tab w:val="left" w:pos="8244"/>
      <w:tab w:val="left" w:pos="9160"/>
      <w:tab w:val="left" w:pos="10076"/>
      <w:tab w:val="left" w:pos="10992"/>
      <w:tab w:val="left" w:pos="11908"/>
      <w:tab w:val="left" w:pos="12824"/>
      <w:tab w:val="left" w:pos="13740"/>
      <w:tab w:val="left" w:pos="14656"/>
    </w:tabs>
    <w:jc w:val="left"/>
    <w:rPr>
      <w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体" w:cs="宋体" w:hint="eastAsia"/>
      <w:color w:val="000000"/>
      <w:kern w:val="0"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体" w:cs="宋体"/>
      <w:color w:val="000000"/>
      <w:kern w:val="0"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>}</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:widowControl/>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:spacing w:before="150" w:after="150"/>
    <w:jc w:val="left"/>
    <w:rPr>
      <w:rFonts w:ascii="微软雅黑" w:eastAsia="微软雅黑" w:hAnsi="微软雅黑" w:cs="宋体"/>
      <w:color w:val="333333"/>
      <w:kern w:val="0"/>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="微软雅黑" w:eastAsia="微软雅黑" w:hAnsi="微软雅黑" w:cs="宋体" w:hint="eastAsia"/>
      <w:color w:val="333333"/>
      <w:kern w:val="0"/>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t>每个pipeline后面有一个数值, 这个数值的范围是0-1000, 这个数值确定了他们的运行顺序(即优先级), 数字越小越优先执行.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>

'@

$target = $d.Paragraphs.Last.Range
[void]$target.InsertXML($xml)
